# "option to clear inactive list" (Id=17) is moved from the Active sheet
# to the Inactive sheet (and marked Done), matching the commit message.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# 1) Remove the task row from the Active sheet (row 4: Id 17).
$active.Rows.Item(4).Delete()

# 2) Insert it at the top of the Inactive sheet's data (row 2) with its
#    new Inactive-sheet fields (Status = Done, last-modified = 4/11/2018).
$inactive.Rows.Item(2).Insert()

# The freshly inserted row picked up the header row's (bold) formatting;
# reset the whole row back to the plain/default style used by every other
# data row before filling in values.
$inactive.Range("A2:F2").Style = "Normal"

$inactive.Range("A2").Value = 17
$inactive.Range("B2").Value = "option to clear inactive list"

# Created/last-modified columns hold plain text that looks like dates, so
# force text formatting before assigning, then strip the number format
# back off so the cell ends up like its neighbours (General, unstyled).
$inactive.Range("E2").NumberFormat = "@"
$inactive.Range("F2").NumberFormat = "@"

$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "Feature"
$inactive.Range("E2").Value = "12/18/2017"
$inactive.Range("F2").Value = "4/11/2018"

$inactive.Range("E2").Style = "Normal"
$inactive.Range("F2").Style = "Normal"
